$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated power-flow results for the 380 kV case (rows 2-25, i.e. data rows 0-23)
$newValues = @{
    2 = @{ "B"=0.5802201919648269; "C"=0.1333860406410849; "D"=0.3486400274081376; "F"=0.927735168085718; "G"=0.3735667452777491; "H"=0.5320773829552436; "I"=0.3757125146676756; "J"=0.2863611567612168; "M"=0.3451593160192559; "N"=1.052571756193949; "O"=1.751971742156471 }
    3 = @{ "B"=0.5095561750466686; "C"=0.116996475682754; "D"=0.345621833688341; "F"=0.9265826203569674; "G"=0.3718514101684534; "H"=0.5347920870870624; "I"=0.3802615564006295; "J"=0.2865529949736825; "M"=0.3217656451734001; "N"=1.056540501829922; "O"=1.753764469106414 }
    4 = @{ "B"=0.4660641078759511; "C"=0.1068845514444661; "D"=0.3439100971284859; "F"=0.9264304020022536; "G"=0.3710938275974272; "H"=0.5367170328413877; "I"=0.3832994015100493; "J"=0.2868265586623693; "M"=0.3074844715521934; "N"=1.059367729072591; "O"=1.756023224164437 }
    5 = @{ "B"=0.4483158692900702; "C"=0.1027519054841548; "D"=0.3432482294262087; "F"=0.9265080882385206; "G"=0.3708593643635965; "H"=0.5375664061375929; "I"=0.384598828156836; "J"=0.2869772109626041; "M"=0.3016859582817446; "N"=1.060618205570854; "O"=1.757234776543157 }
    6 = @{ "B"=0.4453673241770844; "C"=0.1020649690534583; "D"=0.3431404849755495; "F"=0.926529427801654; "G"=0.3708249144644782; "H"=0.5377113674385896; "I"=0.384818308578847; "J"=0.2870045926748688; "M"=0.3007244106705684; "N"=1.060831792471284; "O"=1.757453532147935 }
    7 = @{ "B"=0.4658248478868359; "C"=0.1068288651644309; "D"=0.3439010263396227; "F"=0.9264308839416842; "G"=0.3710903649883193; "H"=0.5367282247859251; "I"=0.383316677175344; "J"=0.2868284317994139; "M"=0.3074061844810956; "N"=1.05938419493355; "O"=1.756038385100567 }
    8 = @{ "B"=0.5558775992833205; "C"=0.1277451836682246; "D"=0.3475700618990345; "F"=0.9272225092255013; "G"=0.3729138859260601; "H"=0.5329598659901009; "I"=0.3772302005196444; "J"=0.2863949722218848; "M"=0.3370762608777298; "N"=1.053859263064027; "O"=1.752349463586157 }
    9 = @{ "B"=0.731596043920149; "C"=0.1683653357121955; "D"=0.3558831553151123; "F"=0.933181658919672; "G"=0.3788401540185617; "H"=0.5276166446873418; "I"=0.3672388236215269; "J"=0.2867812160231225; "M"=0.3959002104096072; "N"=1.046114872553346; "O"=1.754310790238833 }
    10 = @{ "B"=0.8601095478591105; "C"=0.1979556951998802; "D"=0.36266780400058; "F"=0.9402478439415489; "G"=0.3846345347705835; "H"=0.5249370824489574; "I"=0.3610867876453057; "J"=0.2878195157456958; "M"=0.4394936779871088; "N"=1.042299396933544; "O"=1.761370639524841 }
    11 = @{ "B"=0.9184356192523069; "C"=0.2113598610868621; "D"=0.3659004058900024; "F"=0.9440465601166608; "G"=0.3875849794265775; "H"=0.5239883732576374; "I"=0.3585468981721149; "J"=0.2884559234726538; "M"=0.4594040353567976; "N"=1.040968765463106; "O"=1.765805459973564 }
    12 = @{ "B"=0.9405015509258305; "C"=0.2164272729033598; "D"=0.3671454376103895; "F"=0.9455690460841453; "G"=0.388747576938357; "H"=0.5236679521940317; "I"=0.3576223721003302; "J"=0.2887205172341609; "M"=0.4669546458164078; "N"=1.040522977791539; "O"=1.767660889761402 }
    13 = @{ "B"=0.9357502092170762; "C"=0.215336296075634; "D"=0.3668763693194563; "F"=0.94523741619318; "G"=0.3884951732683675; "H"=0.5237352338306493; "I"=0.3578198263474803; "J"=0.2886624825340434; "M"=0.4653280049305835; "N"=1.04061640479047; "O"=1.767253457011662 }
    14 = @{ "B"=0.9202514234396517; "C"=0.2117769313333326; "D"=0.3660024168641911; "F"=0.9441701327147882; "G"=0.3876797180352298; "H"=0.5239612338672259; "I"=0.3584700893889341; "J"=0.2884772187517513; "M"=0.4600250108695292; "N"=1.0409309268166; "O"=1.765954577920809 }
    15 = @{ "B"=0.9107552075819285; "C"=0.2095956073691525; "D"=0.3654698162349774; "F"=0.9435273286294432; "G"=0.387186134046047; "H"=0.524104722064834; "I"=0.358873250760432; "J"=0.2883668129060908; "M"=0.4567781921517806; "N"=1.041131141700944; "O"=1.765181909700971 }
    16 = @{ "B"=0.8562949565534836; "C"=0.1970785351955442; "D"=0.3624594788117292; "F"=0.9400113440285907; "G"=0.3844480540991668; "H"=0.5250045185310057; "I"=0.3612579874740405; "J"=0.2877812271857181; "M"=0.438194051059277; "N"=1.042394494360707; "O"=1.761105442127388 }
    17 = @{ "B"=0.822849674760505; "C"=0.1893849896229369; "D"=0.3606501077535427; "F"=0.9380040317156926; "G"=0.3828489637696464; "H"=0.5256257107903934; "I"=0.3627872588082752; "J"=0.2874640180985253; "M"=0.4268133254594986; "N"=1.043273153481536; "O"=1.758918079725248 }
    18 = @{ "B"=0.8036001767137009; "C"=0.1849545532813295; "D"=0.3596231760747912; "F"=0.9369044747984248; "G"=0.3819588107651981; "H"=0.5260084416574955; "I"=0.3636912011526121; "J"=0.2872970114950135; "M"=0.4202749393891807; "N"=1.043816668780465; "O"=1.757775096574875 }
    19 = @{ "B"=0.7970804969662026; "C"=0.1834535807198563; "D"=0.3592778432022072; "F"=0.9365416300238962; "G"=0.3816625013830759; "H"=0.5261423974384343; "I"=0.3640014396775832; "J"=0.2872431182240902; "M"=0.4180624567506399; "N"=1.044007247929244; "O"=1.757407870967199 }
    20 = @{ "B"=0.8264113023017217; "C"=0.1902045327471455; "D"=0.3608412939893526; "F"=0.9382120218180461; "G"=0.3830161253255397; "H"=0.5255569513093832; "I"=0.3626219450162971; "J"=0.2874961871450239; "M"=0.4280240494243088; "N"=1.04317567312286; "O"=1.75913901153578 }
    21 = @{ "B"=0.9248043706584781; "C"=0.2128226349746001; "D"=0.366258551139822; "F"=0.9444813403467123; "G"=0.3879180060534964; "H"=0.5238937984223782; "I"=0.3582780791533438; "J"=0.2885309946768686; "M"=0.4615823339875078; "N"=1.040836968620411; "O"=1.766331310761416 }
    22 = @{ "B"=0.9889875332610245; "C"=0.2275554468831729; "D"=0.3699209013817466; "F"=0.9490682753827286; "G"=0.3913859191241897; "H"=0.5230331773675516; "I"=0.3556564125847324; "J"=0.289344862246331; "M"=0.4835784692308494; "N"=1.039647036206205; "O"=1.772058173852372 }
    23 = @{ "B"=0.954743439333015; "C"=0.2196968928165859; "D"=0.3679551231343225; "F"=0.9465753539149802; "G"=0.3895108183741769; "H"=0.5234718047620817; "I"=0.3570357380434288; "J"=0.2888978972223697; "M"=0.4718330196255991; "N"=1.040251197680661; "O"=1.768907683826285 }
    24 = @{ "B"=0.8248011569153277; "C"=0.1898340400621805; "D"=0.3607548172573161; "F"=0.9381178198100315; "G"=0.3829404606836562; "H"=0.5255879577289875; "I"=0.3626966063118324; "J"=0.2874815956714656; "M"=0.4274766669309074; "N"=1.043219624484635; "O"=1.759038771427925 }
    25 = @{ "B"=0.6841587423061242; "C"=0.1574201601866321; "D"=0.3535150173839838; "F"=0.9310976149251928; "G"=0.3769845387514721; "H"=0.5288431898323154; "I"=0.3697332646806331; "J"=0.2865442863320737; "M"=0.3799197262180769; "N"=1.047880119658004; "O"=1.752794408939536 }
}

foreach ($rowNum in $newValues.Keys) {
    $rowData = $newValues[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
